$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) - update column F (report count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1713
$ws1.Range("F6").Value = 3288
$ws1.Range("F8").Value = 2119
$ws1.Range("F9").Value = 2042
$ws1.Range("F10").Value = 1053
$ws1.Range("F13").Value = 1636
$ws1.Range("F14").Value = 359
$ws1.Range("F16").Value = 19
$ws1.Range("F18").Value = 122
$ws1.Range("F19").Value = 1494
$ws1.Range("F20").Value = 562
$ws1.Range("F21").Value = 664
$ws1.Range("F23").Value = 11925
$ws1.Range("F24").Value = 11942
$ws1.Range("F25").Value = 878
$ws1.Range("F27").Value = 150
$ws1.Range("F28").Value = 1871
$ws1.Range("F29").Value = 171
$ws1.Range("F30").Value = 494

# Sheet "全部类型" (sheet4.xml) - same data duplicated, update column F values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1713
$ws4.Range("F8").Value = 3288
$ws4.Range("F10").Value = 2119
$ws4.Range("F11").Value = 2042
$ws4.Range("F12").Value = 1053
$ws4.Range("F15").Value = 1636
$ws4.Range("F16").Value = 359
$ws4.Range("F18").Value = 19
$ws4.Range("F22").Value = 122
$ws4.Range("F23").Value = 1494
$ws4.Range("F24").Value = 562
$ws4.Range("F25").Value = 664
$ws4.Range("F27").Value = 11925
$ws4.Range("F28").Value = 11942
$ws4.Range("F29").Value = 878
$ws4.Range("F31").Value = 150
$ws4.Range("F32").Value = 1871
$ws4.Range("F35").Value = 171
$ws4.Range("F36").Value = 494
